$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together as a record, excluding D (Rödlistade) which
# stays put, and C (Valideringsstatus) which is identical across these rows.
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

# Snapshot current values for rows 2-5 before overwriting anything, since
# the update cyclically permutes the rows (row2<-row5, row3<-row4,
# row4<-row3, row5<-row2).
$rows = @(2, 3, 4, 5)
$snapshot = @{}
foreach ($r in $rows) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value()
    }
    $snapshot[$r] = $rowVals
}

# Mapping: new row r gets the old values that used to live at $source[r]
$source = @{ 2 = 5; 3 = 4; 4 = 3; 5 = 2 }

foreach ($r in $rows) {
    $srcRow = $source[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $snapshot[$srcRow][$c]
    }
}
